$wb = $excel.ActiveWorkbook

# --- Neg_Change sheet: replace the 6 existing data rows with 10 new data rows ---
$wsNeg = $wb.Worksheets.Item("Neg_Change")

$negData = @(
        @("HINDZINC", 509, 516.1, 506.35, 513, 6882125, 15420228, -0.5536949907614855, "HINDZINC"),
        @("GODREJCP", 1106, 1118, 1100.7, 1111, 809694, 1831340, -0.5578680092172944, "GODREJCP"),
        @("INDUSTOWER", 339.9, 346.05, 339.4, 345.15, 3966733, 9282025, -0.5726435772366483, "INDUSTOWER"),
        @("PETRONET", 277, 281.35, 276, 280, 1678829, 3379266, -0.5031971439951752, "PETRONET"),
        @("DIXON", 16610, 16875, 16604, 16770, 214714, 455160, -0.5282669830389314, "DIXON"),
        @("PIIND", 3572.6, 3572.6, 3542.7, 3555, 93641, 195195, -0.520269474115628, "PIIND"),
        @("ANGELONE", 2406.9, 2460.7, 2381, 2445, 1343403, 3116178, -0.5688940105475361, "ANGELONE"),
        @("MCX", 9560, 9617, 9425, 9533, 1015096, 2000223, -0.4925085852927399, "MCX"),
        @("KAYNES", 6945, 7012, 6889, 6965, 258385, 556166, -0.5354174832693837, "KAYNES"),
        @("BANDHANBNK", 163.5, 164.5, 162.6, 163.42, 5380906, 10676188, -0.495989954466894, "BANDHANBNK")
)

$r = 2
foreach ($row in $negData) {
    $c = 1
    foreach ($val in $row) {
        $wsNeg.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# --- Pos_Change sheet: replace the 22 existing data rows with 18 new data rows ---
$wsPos = $wb.Worksheets.Item("Pos_Change")

$posData = @(
        @("TRENT", 4617.2, 4752.1, 4610.1, 4739, 1298534, 913650, 0.4212597821923056, "TRENT"),
        @("GRASIM", 2780, 2824.9, 2774.2, 2815, 459912, 316183, 0.4545753566763551, "GRASIM"),
        @("COALINDIA", 381.25, 385.4, 381.25, 384.8, 6135157, 4288494, 0.4306087405042423, "COALINDIA"),
        @("ITC", 398, 400.85, 397.25, 400.3, 16043867, 10291102, 0.5590037879325266, "ITC"),
        @("CIPLA", 1550, 1564.8, 1546, 1556.1, 1638237, 1056222, 0.5510347256542658, "CIPLA"),
        @("KOTAKBANK", 2155, 2165.9, 2141.2, 2148.2, 2308107, 1604563, 0.4384645539003454, "KOTAKBANK"),
        @("MAZDOCK", 2800, 2870, 2786.8, 2850.5, 789866, 533596, 0.4802697171642966, "MAZDOCK"),
        @("DMART", 4230, 4275, 4215, 4264.5, 415827, 276492, 0.5039386311358014, "DMART"),
        @("AMBUJACEM", 563.35, 569, 563.35, 567.75, 1377329, 975792, 0.4114985570695394, "AMBUJACEM"),
        @("VBL", 441.6, 446.5, 440.2, 444.4, 2061374, 1462957, 0.4090461988971651, "VBL"),
        @("DIVISLAB", 6580, 6608.5, 6527, 6556, 462127, 321275, 0.4384156874951365, "DIVISLAB"),
        @("360ONE", 1131, 1172, 1117.3, 1161, 2005553, 1375997, 0.4575271603063088, "360ONE"),
        @("PHOENIXLTD", 1626.2, 1681, 1616.1, 1675, 1339239, 877286, 0.5265705824554364, "PHOENIXLTD"),
        @("ASTRAL", 1418.6, 1449, 1415, 1437.9, 558872, 388990, 0.4367258798426695, "ASTRAL"),
        @("NATIONALUM", 223.8, 228, 223.8, 226.03, 12613195, 8629781, 0.4615892338403489, "NATIONALUM"),
        @("PATANJALI", 583, 589.35, 580.1, 584.6, 556159, 383700, 0.4494631222309096, "PATANJALI"),
        @("GLENMARK", 1892, 1913.5, 1876.2, 1889, 658841, 435764, 0.5119215905857298, "GLENMARK"),
        @("ABFRL", 80.67, 83.65000000000001, 80.22, 83.27, 4595235, 2901302, 0.5838526978577204, "ABFRL")
)

$r = 2
foreach ($row in $posData) {
    $c = 1
    foreach ($val in $row) {
        $wsPos.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# Pos_Change shrank from 23 rows (22 data rows) to 19 rows (18 data rows);
# remove the now-stale trailing rows 20:23 so the used range / dimension shrinks to A1:I19.
$wsPos.Range("A20:I23").Delete()
